$wb = $excel.ActiveWorkbook

# --- Summary sheet: recompute top-level stats after trade #24 closes ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99    # Current Capital
$summary.Range("B4").Value = -0.01      # Total P&L $
$summary.Range("B5").Value = -0.01      # Total P&L %
$summary.Range("B6").Value = 24         # Total Trades
$summary.Range("B7").Value = 7          # Winning Trades
$summary.Range("B9").Value = 29.17      # Win Rate %

# --- Strategy Status sheet: MarketMaking row (row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999   # Capital
$status.Range("D4").Value = 24                  # Trades
$status.Range("E4").Value = -0.01               # P&L $
$status.Range("F4").Value = -0.01               # P&L %
$status.Range("G4").Value = 29.17               # Win Rate %

# --- Append trade #24 to the trade logs (All Trades + MarketMaking) ---
function Add-TradeRow($ws) {
    $r = 25
    $ws.Cells.Item($r, 1).Value = 24
    # B25 looks like an ISO date ("2026-02-17"); Excel would normally
    # auto-convert a bare date-looking string to a date serial on entry.
    # Force it to stay plain text (matching the rest of the Date column),
    # then drop the temporary "@" number format so no stray style sticks.
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = "2026-02-17"
    $cellB.ClearFormats()
    $ws.Cells.Item($r, 3).Value = "15:19:10"
    $ws.Cells.Item($r, 4).Value = "MarketMaking"
    $ws.Cells.Item($r, 5).Value = "DOWN"
    $ws.Cells.Item($r, 6).Value = 0.08
    $ws.Cells.Item($r, 7).Value = 0.27
    $ws.Cells.Item($r, 8).Value = "CLOSED"
    $ws.Cells.Item($r, 9).Value = 237.5
    $ws.Cells.Item($r, 10).Value = 0.19
    $ws.Cells.Item($r, 11).Value = 99.98999999999999
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0.6
    $ws.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($r, 16).Value = "early_exit"
    $ws.Cells.Item($r, 17).Value = 0.15
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$mm = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $mm
